$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - (Intercept)
$ws.Range("B2").Value = 94.756795
$ws.Range("D2").Value = 3.895502
$ws.Range("E2").Value = 0.049244

# Row 3 - household_group_collapsed
$ws.Range("B3").Value = 451.817692
$ws.Range("D3").Value = 9.287233000000001
$ws.Range("E3").Value = 0.000119

# Row 4 - Residuals
$ws.Range("B4").Value = 8075.788968
$ws.Range("C4").Value = 332

# Row 5 - SM-Control
$ws.Range("G5").Value = -1.261358
$ws.Range("H5").Value = -3.049523
$ws.Range("I5").Value = 0.526806
$ws.Range("J5").Value = 0.221988

# Row 6 - SM + Traps-Control
$ws.Range("G6").Value = 1.31937
$ws.Range("H6").Value = -0.558809
$ws.Range("I6").Value = 3.197549
$ws.Range("J6").Value = 0.224739

# Row 7 - SM + Traps-SM
$ws.Range("G7").Value = 2.580728
$ws.Range("H7").Value = 1.168463
$ws.Range("I7").Value = 3.992993
$ws.Range("J7").Value = 0.0000660000000000001
